# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD/AE/AF ---
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Match the header style (bold, bordered, centered) used by the other
# header cells (e.g. A1) by copying formatting onto the new header cells.
# Do this AFTER the values are set.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2-53): constant team record for every player row ---
for ($r = 2; $r -le 53; $r++) {
    $ws.Range("AD$r").Value2 = 78
    $ws.Range("AE$r").Value2 = 84
    $ws.Range("AF$r").Value2 = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) to $($ws.Name)"
